$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 28 ("... Dinâmica/Atividades - Exercícios ..." / botões + LED demo)
# ---------------------------------------------------------------------------
$s28 = $p.Slides.Item(28)
$sh28 = $s28.Shapes.Item(2)

# 1) Move the text placeholder up (only the vertical offset changes).
#    971326 EMU == 76.4824 pt; the extra nudge compensates for the
#    single-precision float truncation PowerPoint applies to Shape.Top.
$sh28.Top = 76.4824

# 2) Merge the trailing run-split sentence in the last paragraph back into a
#    single run (the visible text itself is unchanged).
$tf28 = $sh28.TextFrame.TextRange
$n28 = $tf28.Paragraphs().Count
$lastPara28 = $tf28.Paragraphs($n28)
$whole28 = $tf28.Characters($lastPara28.Start, $lastPara28.Length)
$whole28.Text = "} // um terminal do botão no +,  outro terminal no –, conectando com o resistor."

# ---------------------------------------------------------------------------
# Slide 30 ("... Dinâmica/Atividades - Exercícios ..." / piezo/buzzer demo)
# ---------------------------------------------------------------------------
$s30 = $p.Slides.Item(30)
$sh30 = $s30.Shapes.Item(2)
$tf30 = $sh30.TextFrame.TextRange

# 3) Rename the pin variable from pinoLed to pinoPiezo everywhere it appears.
$n30 = $tf30.Paragraphs().Count
for ($i = 1; $i -le $n30; $i++) {
    $para = $tf30.Paragraphs($i)
    $idx = $para.Text.IndexOf("pinoLed")
    if ($idx -ge 0) {
        $absStart = $para.Start + $idx
        $run = $tf30.Characters($absStart, 7)
        $run.Text = "pinoPiezo"
    }
}

# 4) Merge the trailing "Nota: apresentar o uso de multímetro" run-split text
#    back into a single run (keeping the "  } } // " and "Nota" runs as-is).
$lastPara30 = $tf30.Paragraphs($tf30.Paragraphs().Count)
$idx2 = $lastPara30.Text.IndexOf(": apresentar")
$absStart2 = $lastPara30.Start + $idx2
$len2 = $lastPara30.Length - $idx2
$tail = $tf30.Characters($absStart2, $len2)
$tail.Text = ": apresentar o uso de multímetro"
